$d = $word.ActiveDocument

$replacements = @(
    @{Old="2025-04-27 Sunday"; New="2025-04-28 Monday"},
    @{Old="110×7=770"; New="347×8=2776"},
    @{Old="324×2=648"; New="294×9=2646"},
    @{Old="187×8=1496"; New="348×9=3132"},
    @{Old="159×7=1113"; New="418×4=1672"},
    @{Old="903×2=1806"; New="267×9=2403"},
    @{Old="578×5=2890"; New="601×9=5409"},
    @{Old="828×6=4968"; New="158×9=1422"},
    @{Old="932×9=8388"; New="178×2=356"},
    @{Old="502×8=4016"; New="529×3=1587"},
    @{Old="231×3=693"; New="987×3=2961"},
    @{Old="199×7=1393"; New="316×4=1264"},
    @{Old="852×4=3408"; New="398×7=2786"},
    @{Old="840×6=5040"; New="442×3=1326"},
    @{Old="414×2=828"; New="879×5=4395"},
    @{Old="625×3=1875"; New="376×4=1504"},
    @{Old="596×3=1788"; New="386×3=1158"},
    @{Old="589×9=5301"; New="790×5=3950"},
    @{Old="745×3=2235"; New="625×4=2500"},
    @{Old="715×6=4290"; New="542×6=3252"},
    @{Old="400×2=800"; New="704×4=2816"},
    @{Old="320×2=640"; New="251×7=1757"},
    @{Old="827×6=4962"; New="562×8=4496"},
    @{Old="303×7=2121"; New="367×8=2936"},
    @{Old="873×6=5238"; New="716×3=2148"},
    @{Old="269×8=2152"; New="445×6=2670"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
